$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-01-02 Thursday"; new = "2025-01-03 Friday"},
    @{old = "327×4=1308"; new = "380×7=2660"},
    @{old = "357×4=1428"; new = "603×2=1206"},
    @{old = "988×4=3952"; new = "514×4=2056"},
    @{old = "142×2=284"; new = "991×9=8919"},
    @{old = "578×8=4624"; new = "215×6=1290"},
    @{old = "753×5=3765"; new = "210×9=1890"},
    @{old = "131×6=786"; new = "707×2=1414"},
    @{old = "214×5=1070"; new = "658×3=1974"},
    @{old = "249×5=1245"; new = "800×8=6400"},
    @{old = "798×4=3192"; new = "120×7=840"},
    @{old = "470×7=3290"; new = "434×7=3038"},
    @{old = "867×6=5202"; new = "443×4=1772"},
    @{old = "954×4=3816"; new = "328×3=984"},
    @{old = "240×7=1680"; new = "974×6=5844"},
    @{old = "389×3=1167"; new = "426×5=2130"},
    @{old = "342×7=2394"; new = "847×4=3388"},
    @{old = "452×2=904"; new = "211×5=1055"},
    @{old = "281×4=1124"; new = "662×5=3310"},
    @{old = "467×7=3269"; new = "160×8=1280"},
    @{old = "440×3=1320"; new = "439×2=878"},
    @{old = "668×4=2672"; new = "937×8=7496"},
    @{old = "512×6=3072"; new = "304×4=1216"},
    @{old = "677×9=6093"; new = "997×9=8973"},
    @{old = "778×8=6224"; new = "816×8=6528"},
    @{old = "806×7=5642"; new = "743×5=3715"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
